$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("test")
$ws.Activate()
$ws.Range("A20").Value = "5/25/20"
$ws.Range("B20").Value = "2.99"
$ws.Range("C20").Value = "40:42"
$ws.Range("D20").Value = "13:36"

$ws.Range("A21").Value = "2"
$ws.Range("B21").Value = "123"
$ws.Range("C21").Value = "12:12"
$ws.Range("D21").Value = "0:5"

$ws.Range("A22").Value = "4/22/20"
$ws.Range("B22").Value = "32"
$ws.Range("C22").Value = "312"

$ws.Range("A23").Value = "33123"
$ws.Range("B23").Value = "32"
$ws.Range("C23").Value = "32:32"
$ws.Range("D23").Value = "1:1"

$ws.Range("A24").Value = "05/26/2020"
$ws.Range("B24").Value = "1.58"
$ws.Range("C24").Value = "24:38"
$ws.Range("D24").Value = "15:35"
